$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 165-166, shifting the existing rows 165:190 down to 167:192.
$ws.Rows("165:166").Insert()

# Row 165 - new weekly entry (Primera)
$ws.Cells.Item(165, 1).Value = 1
$ws.Cells.Item(165, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(165, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(165, 4).Value = 44522
$ws.Cells.Item(165, 5).Value = 15
$ws.Cells.Item(165, 6).Value = 100114014
$ws.Cells.Item(165, 7).Value = "Betarraga"
$ws.Cells.Item(165, 8).Value = "Sin especificar"
$ws.Cells.Item(165, 9).Value = "Primera"
$ws.Cells.Item(165, 10).Value = 800
$ws.Cells.Item(165, 11).Value = 300
$ws.Cells.Item(165, 12).Value = 350
$ws.Cells.Item(165, 13).Value = 325
$ws.Cells.Item(165, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(165, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(165, 16).Value = 81
$ws.Cells.Item(165, 17).Value = 4
$ws.Cells.Item(165, 18).Value = "Hortaliza"

# Row 166 - new weekly entry (Segunda)
$ws.Cells.Item(166, 1).Value = 1
$ws.Cells.Item(166, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(166, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(166, 4).Value = 44522
$ws.Cells.Item(166, 5).Value = 15
$ws.Cells.Item(166, 6).Value = 100114014
$ws.Cells.Item(166, 7).Value = "Betarraga"
$ws.Cells.Item(166, 8).Value = "Sin especificar"
$ws.Cells.Item(166, 9).Value = "Segunda"
$ws.Cells.Item(166, 10).Value = 1000
$ws.Cells.Item(166, 11).Value = 300
$ws.Cells.Item(166, 12).Value = 350
$ws.Cells.Item(166, 13).Value = 325
$ws.Cells.Item(166, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(166, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(166, 16).Value = 65
$ws.Cells.Item(166, 17).Value = 5
$ws.Cells.Item(166, 18).Value = "Hortaliza"

# Make sure the date cells keep the original date number format (style carried by Insert,
# but set explicitly to be safe).
$ws.Range("D165:D166").NumberFormat = $ws.Range("D164").NumberFormat
